$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.170.72'
$ws.Range('E2').Value = '  +0.67%  '

$ws.Range('D3').Value = '4.020.06'
$ws.Range('E3').Value = '  -0.17%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '''530.75'
$ws.Range('E5').Value = '  +1.07%  '

$ws.Range('D6').Value = '''151.22'
$ws.Range('E6').Value = '  +1.36%  '

$ws.Range('D7').Value = '''0.696'
$ws.Range('E7').Value = '  +11.28%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  +0.75%  '

$ws.Range('E11').Value = '  -4.34%  '

$ws.Range('D12').Value = '''47.82'
$ws.Range('E12').Value = '  +3.78%  '

$ws.Range('D13').Value = '''10.66'
$ws.Range('E13').Value = '  -1.11%  '

$ws.Range('D14').Value = '4.662.14'
$ws.Range('E14').Value = '  -0.17%  '

$ws.Range('D15').Value = '4.008.46'
$ws.Range('E15').Value = '  -0.66%  '

$ws.Range('D16').Value = '''14.13'
$ws.Range('E16').Value = '  -1.23%  '

$ws.Range('D17').Value = '''20.57'
$ws.Range('E17').Value = '  -4.19%  '

$ws.Range('D18').Value = '''0.133'
$ws.Range('E18').Value = '  -0.64%  '

$ws.Range('E19').Value = '  -2.51%  '

$ws.Range('D20').Value = '72.011.53'
$ws.Range('E20').Value = '  +0.55%  '

$ws.Range('D21').Value = '''427.23'
$ws.Range('E21').Value = '  -3.17%  '

$ws.Range('D22').Value = '''98.21'
$ws.Range('E22').Value = '  +3.57%  '

$ws.Range('E23').Value = '  -3.86%  '

$ws.Range('D24').Value = '''4.22'
$ws.Range('E24').Value = '  +3.73%  '

$ws.Range('D25').Value = '''14.34'
$ws.Range('E25').Value = '  -0.38%  '

$ws.Range('D26').Value = '''11.24'
$ws.Range('E26').Value = '  -7.98%  '

$ws.Range('D27').Value = '''10.76'
$ws.Range('E27').Value = '  -4.16%  '

$ws.Range('E28').Value = '  +0.97%  '

$ws.Range('D29').Value = '''36.68'
$ws.Range('E29').Value = '  -1.04%  '

$ws.Range('D30').Value = '''3.55'
$ws.Range('E30').Value = '  +21.80%  '

$ws.Range('D31').Value = '''13.39'
$ws.Range('E31').Value = '  -1.74%  '

$ws.Range('D32').Value = '''7.17'
$ws.Range('E32').Value = '  +3.26%  '

$ws.Range('E33').Value = '  -2.18%  '

$ws.Range('D34').Value = '''677.07'
$ws.Range('E34').Value = '  -3.25%  '

$ws.Range('D35').Value = '''44.52'
$ws.Range('E35').Value = '  +8.43%  '

$ws.Range('D36').Value = '''65.71'
$ws.Range('E36').Value = '  -3.13%  '

$ws.Range('D37').Value = '''0.448'
$ws.Range('E37').Value = '  +0.83%  '

$ws.Range('E38').Value = '  -3.82%  '

$ws.Range('D39').Value = '0.0₃0828'
$ws.Range('E39').Value = '  -8.85%  '

$ws.Range('D40').Value = '''3.38'
$ws.Range('E40').Value = '  -5.60%  '

$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  -0.11%  '

$ws.Range('E42').Value = '  -0.14%  '

$ws.Range('D43').Value = '''0.0487'
$ws.Range('E43').Value = '  -0.84%  '

$ws.Range('D44').Value = '''3.19'
$ws.Range('E44').Value = '  +1.84%  '

$ws.Range('E45').Value = '  +2.46%  '

$ws.Range('D46').Value = '''3.46'
$ws.Range('E46').Value = '  -1.43%  '

$ws.Range('D47').Value = '''9.71'
$ws.Range('E47').Value = '  +5.62%  '

$ws.Range('E48').Value = '  -7.20%  '

$ws.Range('E49').Value = '  -6.17%  '

# Row 50: Monero -> FLOKI
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = '''0.000273'
$ws.Range('E50').Value = '  -1.56%  '

# Row 51: FLOKI -> Monero
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''146.24'
$ws.Range('E51').Value = '  +1.82%  '
